$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 742 ("「開⇔閉」" entry), shifting all subsequent rows up by one.
$ws.Rows.Item(742).Delete()
